# Adds 3 new shortage items (CETAL, MICONAZ, شامبو جونسون 200مللى) to the
# "نواقص الأصناف" (shortage items) report, keeping the existing alphabetical
# ordering, updates the printed total and refreshes the "printed at" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert three new data rows right after the current
#        last item row (row 12), before the totals row (old row 13). ---
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(15).Insert()

# Copy the formatting of an existing item row onto each of the new rows so
# they look like the rest of the table (borders/fonts/number formats).
$ws.Range("A12:Q12").Copy()
$ws.Range("A13:Q15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights as used by the generator for this report.
$ws.Rows.Item(13).RowHeight = 24.75
$ws.Rows.Item(14).RowHeight = 25.5
$ws.Rows.Item(15).RowHeight = 24.75
$ws.Rows.Item(16).RowHeight = 25.5

# Re-create the merged cells for the three new rows (same layout as every
# other item row: A:B, C:G, H:K, L:M, N:O).
$ws.Range("A13:B13").Merge()
$ws.Range("C13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()

$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

$ws.Range("A15:B15").Merge()
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

# --- 2. Re-write the full item list (9 rows) in alphabetical order, now
#        including the three newly stocked-out items. ---
$items = @(
    @{ Row=7;  Name="AMARYL 2 MG 30 TABS";              Balance="0:0"; Limit=1; Price="78.00";  Sale="51.4800"; Trans="0:2" },
    @{ Row=8;  Name="CATAFLAM 75MG/3ML 6 AMP.";          Balance="1:0"; Limit=1; Price="120.00"; Sale="19.2000"; Trans="0:1" },
    @{ Row=9;  Name="CETAL 100MG/ML ORAL DROPS 15 ML";   Balance="4:0"; Limit=1; Price="23.00";  Sale="23.0000"; Trans="1:0" },
    @{ Row=10; Name="DIMRA 20 F.C.TAB.";                 Balance="0:1"; Limit=1; Price="70.00";  Sale="70.0000"; Trans="1:0" },
    @{ Row=11; Name="MICONAZ 2% ORAL GEL 20 GM";         Balance="1:0"; Limit=1; Price="23.00";  Sale="23.0000"; Trans="1:0" },
    @{ Row=12; Name="SILDEN 100 MG 8F.C. TABS";          Balance="2:7"; Limit=1; Price="66.00";  Sale="66.0000"; Trans="1:0" },
    @{ Row=13; Name="سرنجات 3 سم";                       Balance="0:0"; Limit=0; Price="2.00";   Sale="2.0000";  Trans="1:0" },
    @{ Row=14; Name="شامبو جونسون 200مللى";              Balance="1:0"; Limit=0; Price="50.00";  Sale="50.0000"; Trans="1:0" },
    @{ Row=15; Name="كالونا ";                           Balance="0:0"; Limit=0; Price="15.00";  Sale="15.0000"; Trans="1:0" }
)

$seq = 1
foreach ($item in $items) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $seq          # A - م (sequence number)
    $ws.Cells.Item($r, 3).Value = $item.Name    # C - الاسم
    $ws.Cells.Item($r, 8).Value = $item.Balance # H - الرصيد الحالي
    $ws.Cells.Item($r, 12).Value = $item.Limit  # L - حد الطلب
    $ws.Cells.Item($r, 14).Value = $item.Price  # N - السعر
    $ws.Cells.Item($r, 16).Value = $item.Sale   # P - سعر البيع
    $ws.Cells.Item($r, 17).Value = $item.Trans  # Q - عدد التعاملات
    $seq = $seq + 1
}

# --- 3. Update the printed total (sum of the "سعر البيع" column). ---
$ws.Cells.Item(16, 16).Value = 319.68

# --- 4. Refresh the "printed at" timestamp shown in the footer. ---
$ws.Cells.Item(17, 1).Value = "Monday, 4 August, 2025 10:24 AM"
